$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing columns A-D to B-E.
$ws.Columns("A").Insert()

# Give the newly-inserted column A (rows 2-15) the same formatting as the
# header row's style (bold font, thin border, centered/top alignment) by
# copying the format from B1, which already carries that style.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2:A15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update header text to add spaces around the "=" sign.
$ws.Range("D1").Value = "Treatment at T2 (n = 5080)"
$ws.Range("E1").Value = "Control at T1 (n = 745)"

# Update the section-label cells (now in column B after the shift) to add
# spaces around the "=" sign.
$ws.Range("B3").Value = "Gender (P = 0.006)"
$ws.Range("B10").Value = "Interested in News (P = 0.000)"
